$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# This edit inserts a brand-new "零字节" entry as row 2 of the table
# and pushes the previously-existing company rows down by one row
# (row 2 -> 3, row 3 -> 4, ... row 9 -> 10). The sequential index in
# column A is left untouched for the rows that already existed and a
# new value is appended for the newly created row 10.
# ------------------------------------------------------------------

# Start from a clean slate for the data columns (B:P) across every
# row that is affected, so no stale values are left behind once the
# final content is written back out below.
$ws.Range("B2:P10").ClearContents()

# Column A keeps its existing sequential numbering for rows 2-9; only
# row 10 is new and needs both a value and the same bold/border/
# centered style (style index 1, same as the header row) used by the
# rest of column A.
$ws.Range("A10").Value = 8
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)

# Row 2
$ws.Range("B2").Value = "零字节"
$ws.Range("C2").Value = "建邺"
$ws.Range("E2").Value = "Go/Rust/JS/TS/产品/运营"
$ws.Range("F2").Value = "9：30-6：30"
$ws.Range("G2").Value = "1.5h"
$ws.Range("H2").Value = "不加班"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "8%"
$ws.Range("I2").ClearFormats()
$ws.Range("J2").Value = "13薪，每年调薪一次"
$ws.Range("K2").Value = "应届生八折，有工作经验的不打折"
$ws.Range("L2").Value = "macbook pro（入职满三年电脑转赠给员工），每人配一个显示器（24-32寸）"
$ws.Range("M2").Value = "入职转正就享受年假"
$ws.Range("N2").Value = "飞书打卡"
$ws.Range("O2").Value = "节日红包、年度旅游（21年三亚一周）"

# Row 3
$ws.Range("B3").Value = "南京力方科技有限公司(力方智充)"
$ws.Range("C3").Value = "雨花台区软件谷科创城"
$ws.Range("D3").Value = "技术部"
$ws.Range("E3").Value = "Java"
$ws.Range("F3").Value = "9:00-18:00"
$ws.Range("G3").Value = "1.5h"
$ws.Range("H3").Value = "124固定加班到9点"
$ws.Range("I3").Value = "最低，双边合计512"
$ws.Range("J3").Value = "无"
$ws.Range("K3").Value = "三个月，打八折"
$ws.Range("L3").Value = "网吧工位，自带电脑"
$ws.Range("M3").Value = "法定年假"
$ws.Range("N3").Value = "严格打卡，迟打卡扣30，不打卡半天工资"

# Row 4
$ws.Range("B4").Value = "硅基智能"
$ws.Range("C4").Value = "软件大道"
$ws.Range("D4").Value = "创新产品事业群"
$ws.Range("E4").Value = "Java"
$ws.Range("F4").Value = "9:00-18:30"
$ws.Range("G4").Value = "1.5h"
$ws.Range("H4").Value = "没事到点走，部门氛围卷"
$ws.Range("I4").Value = "基数5500，比例10%"
$ws.Range("J4").Value = "13薪还是根据公司业绩提供，是否折扣，折扣比例。"
$ws.Range("K4").Value = "不打折"
$ws.Range("L4").Value = "网吧工位"
$ws.Range("M4").Value = "满一年才有正常年假，年假次年一月发放（不满一年打折）"
$ws.Range("N4").Value = "是否严格打卡，使用的软件或者方式（比如钉钉或人脸识别）。"

# Row 5
$ws.Range("B5").Value = "百家云"
$ws.Range("C5").Value = "雨花台软件谷科创城"
$ws.Range("E5").Value = "Java"
$ws.Range("F5").Value = "9:00-18:30"
$ws.Range("G5").Value = "1.5h"
$ws.Range("H5").Value = "周1,2,4正常加班，不想加班也行"
$ws.Range("K5").Value = "6个月不打折。"
$ws.Range("L5").Value = "mac笔记本+小米曲面屏显示器"
$ws.Range("M5").Value = "年假次年一月发放，每满一年+1天"
$ws.Range("N5").Value = "每个月有4次迟到补卡机会，早上9.15之前打卡不算迟到"

# Row 6
$ws.Range("B6").Value = "创维南京分公司"
$ws.Range("C6").Value = "雨花云密城"
$ws.Range("D6").Value = "web后台"
$ws.Range("E6").Value = "Java"
$ws.Range("F6").Value = "09:30"
$ws.Range("G6").Value = "1.5h"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "995"
$ws.Range("H6").ClearFormats()
$ws.Range("I6").Value = "工资八折的10%"
$ws.Range("J6").Value = "1个月工资"
$ws.Range("K6").Value = "不打折"
$ws.Range("L6").Value = "Windows电脑+dell显示器"
$ws.Range("M6").Value = "法定年假"
$ws.Range("N6").Value = "弹性打卡"

# Row 7
$ws.Range("B7").Value = "新视云"
$ws.Range("C7").Value = "雨花台"
$ws.Range("E7").Value = "Java"
$ws.Range("F7").Value = "9:00-17:30"
$ws.Range("G7").Value = "1h"
$ws.Range("H7").Value = "看部门，业务部门偶尔加班，技术支持部门基本不加班"
$ws.Range("I7").Value = "基数5k，比例8%"
$ws.Range("J7").Value = "固定13薪"
$ws.Range("K7").Value = "3年合同，试用期总共6个月，前三个月8折，后三个月全薪"
$ws.Range("L7").Value = "配笔记本+显示器"
$ws.Range("M7").Value = "5天年假+5天带薪病假（入职自动折算当年年假）"
$ws.Range("N7").Value = "不打卡"

# Row 8
$ws.Range("B8").Value = "华为"
$ws.Range("C8").Value = "华为南研所"
$ws.Range("E8").Value = "Java"
$ws.Range("F8").Value = "9:00"
$ws.Range("G8").Value = "12:00-13:40"
$ws.Range("H8").Value = "看部门情况。好部门：124加班8：30，35正常下班,差部门：天天11点以后"
$ws.Range("I8").Value = "基础工资的5%"
$ws.Range("J8").Value = "看部门盈利情况和个人绩效定"
$ws.Range("K8").Value = "试用期6个月，100%工资不打折"
$ws.Range("L8").Value = "配win台式机+双屏"
$ws.Range("M8").Value = "没签奋斗协议的5天，但一般不给休，第二年可以换成钱。签了的自愿放弃年假了"
$ws.Range("N8").Value = "必须按时打卡"

# Row 9
$ws.Range("B9").Value = "满帮"
$ws.Range("C9").Value = "雨花区万博科技园"
$ws.Range("E9").Value = "Java"
$ws.Range("F9").Value = "9:00-18:30"
$ws.Range("G9").Value = "1.5h"
$ws.Range("H9").Value = "看部门，不强制， 周五基本不加，还有每月一天奋斗日（年底算工资）， 据说要取消了"
$ws.Range("I9").Value = "全额8%"
$ws.Range("J9").Value = "上下半年绩效"
$ws.Range("L9").Value = "联想"

# Row 10
$ws.Range("B10").Value = "A示例xxx公司"
$ws.Range("C10").Value = "xx区"
$ws.Range("D10").Value = "xxx事业部"
$ws.Range("E10").Value = "Java"
$ws.Range("F10").Value = "9:00-18:30"
$ws.Range("G10").Value = "1.5h"
$ws.Range("H10").Value = "135 加班，24 正常；大小周等等"
$ws.Range("I10").Value = "基数 xxxx，比例 xx%"
$ws.Range("J10").Value = "13薪还是根据公司业绩提供，是否折扣，折扣比例。"
$ws.Range("K10").Value = "是否打折，比如 xx%。"
$ws.Range("L10").Value = "工位大小，环境，是否提供设备，设备型号种类。"
$ws.Range("M10").Value = "是否有入职就有，是否有前置条件才有。"
$ws.Range("N10").Value = "是否严格打卡，使用的软件或者方式（比如钉钉或人脸识别）。"
